# New weekly price report: insert two new rows (new week's data) right
# above the previous first data rows, pushing all existing data rows
# down by two (the sheet already kept rows in reverse-chronological
# order). The two newly reported "Poroto verde" entries (Magnum /
# Sin especificar, both "Primera") go into the freed-up rows 12 and 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 12:124 down to 14:126, carrying formatting along.
$ws.Rows("12:13").Insert()

# New row 12 — Poroto verde, Magnum, Primera
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44552
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100112031
$ws.Range("G12").Value = "Poroto verde"
$ws.Range("H12").Value = "Magnum"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("N12").Value = "$/malla 25 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 620
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"

# New row 13 — Poroto verde, Sin especificar, Primera
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44552
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112031
$ws.Range("G13").Value = "Poroto verde"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 28000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 29000
$ws.Range("N13").Value = "$/malla 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 1160
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
